# "Registration Page test Data 2"
# Updates test data on the "reg" worksheet:
#   - Refreshes several "Username" sample values (column C, rows 2-9)
#   - Fills in the previously-empty "Email" value for row 9 (was ""),
#     matching the formatting used by the other Email cells in that column
#   - Moves the active selection from C9 to A9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reg")

# Refresh Username sample values (column C)
$ws.Range("C2").Value = "akhilbingi32131"
$ws.Range("C3").Value = "akhilbingi64234442"
$ws.Range("C4").Value = "akhilbingi3244553"
$ws.Range("C5").Value = "akhilbingi3215664"
$ws.Range("C6").Value = "akhilbingi32127775"
$ws.Range("C7").Value = "nffmf886"
$ws.Range("C8").Value = "akhilhdiw5997"
$ws.Range("C9").Value = "akhilhdingi9008"

# Row 9, Email column (A9): was an empty quoted string with the
# Calibri-based style; give it the same formatting as the other
# populated Email cells (e.g. A2) before setting its value.
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "akhilbingi@gmail.com"

# Move the selection to A9 (was C9)
$ws.Range("A9").Select()
